$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I (I0) and J (IF), copying the
# formatting from the existing header cell H1 so the new headers
# match the style of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J55
$data = @(
    @(7, 7),
    @(5, 6),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(1, 1),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(5, 6),
    @(7, 8),
    @(6, 7),
    @(6, 6),
    @(10, 10),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(6, 7),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(6, 8),
    @(4, 5),
    @(5, 7),
    @(8, 9),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(4, 4)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = 2 + $idx
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}

Write-Output "I0 and IF columns added"